$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Create the new "2022-Q4" sheet by copying the existing "2022-Q3"
#    sheet (so header row / column styles match the other quarter
#    sheets exactly), placing it right before "2022-Q3" (i.e. right
#    after "总计").
# ------------------------------------------------------------------
$src = $wb.Worksheets.Item("2022-Q3")
$summary = $wb.Worksheets.Item("总计")
$src.Copy($null, $summary)
$q4 = $wb.ActiveSheet
$q4.Name = "2022-Q4"

# The source sheet had 4 data rows (rows 2-5); the new sheet only
# needs 2, so drop the extra two rows.
$q4.Rows.Item(4).Delete()
$q4.Rows.Item(4).Delete()

# Make sure the fund-code / numeric-looking text columns stay text
# (matches the "looks numeric but stored as text" convention used by
# every other quarter sheet), then fill in the 2022-Q4 fund data.
$q4.Range("B2:B3").NumberFormat = "@"
$q4.Range("D2:G3").NumberFormat = "@"

$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "009330"
$q4.Range("C2").Value = "鹏华成长价值混合A"
$q4.Range("D2").Value = "7.36"
$q4.Range("E2").Value = "71.66"
$q4.Range("F2").Value = "2.45"
$q4.Range("G2").Value = "0.1803"
$q4.Range("H2").Value = 10

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "009331"
$q4.Range("C3").Value = "鹏华成长价值混合C"
$q4.Range("D3").Value = "3.20"
$q4.Range("E3").Value = "71.66"
$q4.Range("F3").Value = "2.45"
$q4.Range("G3").Value = "0.0784"
$q4.Range("H3").Value = 10

# ------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new row right after
#    the header for the 2022-Q4 totals, pushing every other quarter
#    down by one row.
# ------------------------------------------------------------------
$summary.Rows.Item(2).Insert()

# The freshly inserted row inherits odd blended formatting from the
# Insert() call - clear it so it matches the plain (unstyled) data
# rows used throughout this sheet, then restore column A's index
# style (shared with every other data row) by copying it over.
$summary.Range("B2:D2").ClearFormats()
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.26

# Column A is a plain 0-based row counter that gets renumbered after
# the new row is prepended (not the old values shifted down), so
# rewrite it for every remaining data row.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
$summary.Range("A8").Value = 6
$summary.Range("A9").Value = 7
$summary.Range("A10").Value = 8
